$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1039.44
$ws.Range("I19").Value = 488.875
$ws.Range("J19").Value = 1298.5294
$ws.Range("K19").Value = 488.875
$ws.Range("L19").Value = 1298.5294
$ws.Range("M19").Value = -313.875
$ws.Range("N19").Value = -1648.5294

$ws.Range("H106").Value = 2499.4167
$ws.Range("I106").Value = 2599.1428
$ws.Range("J106").Value = 2359.8
$ws.Range("K106").Value = 2599.1428
$ws.Range("L106").Value = 2359.8
$ws.Range("M106").Value = -1968.1428
$ws.Range("N106").Value = -3621.8

$ws.Range("H113").Value = 78685.46000000001
$ws.Range("I113").Value = 251426.25
$ws.Range("J113").Value = 1911.7778
$ws.Range("K113").Value = 251426.25
$ws.Range("L113").Value = 1911.7778
$ws.Range("M113").Value = -248172.25
$ws.Range("N113").Value = -8419.7778

$ws.Range("H132").Value = 4634274.5
$ws.Range("I132").Value = 4812265.5
$ws.Range("J132").Value = 6503
$ws.Range("K132").Value = 14436796.5
$ws.Range("L132").Value = 19509
$ws.Range("M132").Value = -14434266.5
$ws.Range("N132").Value = -24569

$ws.Range("H135").Value = 941.86487
$ws.Range("I135").Value = 583.6061
$ws.Range("K135").Value = 5252.4549
$ws.Range("M135").Value = -2717.4549

$ws.Range("H138").Value = 5856.6963
$ws.Range("I138").Value = 1131.6586
$ws.Range("J138").Value = 18771.8
$ws.Range("K138").Value = 3394.9758
$ws.Range("L138").Value = 56315.39999999999
$ws.Range("M138").Value = 1745.0242
$ws.Range("N138").Value = -66595.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 166.33333
$ws.Range("I4").Value = 199.5
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 199.5
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -83.5
$ws.Range("N4").Value = -332

$ws.Range("H61").Value = 1707.4615
$ws.Range("I61").Value = 1105.875
$ws.Range("J61").Value = 2670
$ws.Range("K61").Value = 1105.875
$ws.Range("L61").Value = 2670
$ws.Range("M61").Value = -893.875
$ws.Range("N61").Value = -3094

$ws.Range("H63").Value = 1920.8334
$ws.Range("I63").Value = 1391.6666
$ws.Range("J63").Value = 2450
$ws.Range("K63").Value = 1391.6666
$ws.Range("L63").Value = 2450
$ws.Range("M63").Value = -705.6666
$ws.Range("N63").Value = -3822

$ws.Range("H66").Value = 1920.8334
$ws.Range("I66").Value = 1391.6666
$ws.Range("J66").Value = 2450
$ws.Range("K66").Value = 6958.333000000001
$ws.Range("L66").Value = 12250
$ws.Range("M66").Value = -3526.333000000001
$ws.Range("N66").Value = -19114

$ws.Range("H74").Value = 2673.2
$ws.Range("I74").Value = 3804
$ws.Range("J74").Value = 2188.5715
$ws.Range("K74").Value = 3804
$ws.Range("L74").Value = 2188.5715
$ws.Range("M74").Value = -2930
$ws.Range("N74").Value = -3936.5715

$ws.Range("H77").Value = 2673.2
$ws.Range("I77").Value = 3804
$ws.Range("J77").Value = 2188.5715
$ws.Range("K77").Value = 19020
$ws.Range("L77").Value = 10942.8575
$ws.Range("M77").Value = -14652
$ws.Range("N77").Value = -19678.8575

$ws.Range("H122").Value = 1132.037
$ws.Range("I122").Value = 1047.3043
$ws.Range("K122").Value = 3141.9129
$ws.Range("M122").Value = -691.9129000000003

$ws.Range("H132").Value = 3263.147
$ws.Range("I132").Value = 3412.926
$ws.Range("J132").Value = 2685.4285
$ws.Range("K132").Value = 10238.778
$ws.Range("L132").Value = 8056.2855
$ws.Range("M132").Value = -7708.778
$ws.Range("N132").Value = -13116.2855

$ws.Range("H136").Value = 1707.4615
$ws.Range("I136").Value = 1105.875
$ws.Range("J136").Value = 2670
$ws.Range("K136").Value = 3317.625
$ws.Range("L136").Value = 8010
$ws.Range("M136").Value = -767.625
$ws.Range("N136").Value = -13110

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 21097
$ws.Range("J35").Value = 21097
$ws.Range("L35").Value = 21097
$ws.Range("N35").Value = -21717

$ws.Range("H82").Value = 16763.092
$ws.Range("I82").Value = 2796.25
$ws.Range("J82").Value = 24744.143
$ws.Range("K82").Value = 2796.25
$ws.Range("L82").Value = 24744.143
$ws.Range("M82").Value = -2413.25
$ws.Range("N82").Value = -25510.143

$ws.Range("H85").Value = 16763.092
$ws.Range("I85").Value = 2796.25
$ws.Range("J85").Value = 24744.143
$ws.Range("K85").Value = 2796.25
$ws.Range("L85").Value = 24744.143
$ws.Range("M85").Value = -1470.25
$ws.Range("N85").Value = -27396.143

$ws.Range("H134").Value = 2069.4153
$ws.Range("I134").Value = 1813.7966
$ws.Range("J134").Value = 4583
$ws.Range("K134").Value = 5441.3898
$ws.Range("L134").Value = 13749
$ws.Range("M134").Value = -2906.3898
$ws.Range("N134").Value = -18819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 24799.941
$ws.Range("I31").Value = 861.71875
$ws.Range("J31").Value = 46078.36
$ws.Range("K31").Value = 861.71875
$ws.Range("L31").Value = 46078.36
$ws.Range("M31").Value = -566.71875
$ws.Range("N31").Value = -46668.36

$ws.Range("H34").Value = 24799.941
$ws.Range("I34").Value = 861.71875
$ws.Range("J34").Value = 46078.36
$ws.Range("K34").Value = 861.71875
$ws.Range("L34").Value = 46078.36
$ws.Range("M34").Value = -659.71875
$ws.Range("N34").Value = -46482.36

$ws.Range("H119").Value = 39995.5
$ws.Range("J119").Value = 39995.5
$ws.Range("L119").Value = 39995.5
$ws.Range("N119").Value = -49671.5

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H122").Value = 749.5
$ws.Range("I122").Value = 749.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2248.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 201.5
$ws.Range("N122").ClearContents()

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H134").Value = 1025.138
$ws.Range("I134").Value = 888.125
$ws.Range("K134").Value = 2664.375
$ws.Range("M134").Value = -129.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 388.21054
$ws.Range("I18").Value = 321.53333
$ws.Range("K18").Value = 964.5999899999999
$ws.Range("M18").Value = -795.5999899999999

$ws.Range("H122").Value = 453.58334
$ws.Range("I122").Value = 449.22223
$ws.Range("J122").Value = 466.66666
$ws.Range("K122").Value = 4043.00007
$ws.Range("L122").Value = 4199.99994
$ws.Range("M122").Value = -1593.00007
$ws.Range("N122").Value = -9099.99994

$ws.Range("H129").Value = 14175836
$ws.Range("I129").Value = 83352430
$ws.Range("J129").Value = 340517.6
$ws.Range("K129").Value = 250057290
$ws.Range("L129").Value = 1021552.8
$ws.Range("M129").Value = -250052290
$ws.Range("N129").Value = -1031552.8

$ws.Range("H131").Value = 6859.0107
$ws.Range("I131").Value = 848
$ws.Range("J131").Value = 7566.1885
$ws.Range("K131").Value = 2544
$ws.Range("L131").Value = 22698.5655
$ws.Range("M131").Value = 2496
$ws.Range("N131").Value = -32778.5655

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3096.946
$ws.Range("I132").Value = 2051.9048
$ws.Range("K132").Value = 6155.714399999999
$ws.Range("M132").Value = -3625.714399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3513.818
$ws.Range("I132").Value = 3265.3
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 9795.900000000001
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -7265.900000000001
$ws.Range("N132").Value = -23057

$ws.Range("H136").Value = 1623.3448
$ws.Range("I136").Value = 1254.8572
$ws.Range("J136").Value = 2590.625
$ws.Range("K136").Value = 3764.5716
$ws.Range("L136").Value = 7771.875
$ws.Range("M136").Value = -1214.5716
$ws.Range("N136").Value = -12871.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2286.1765
$ws.Range("I132").Value = 1996.4565
$ws.Range("J132").Value = 4951.6
$ws.Range("K132").Value = 5989.3695
$ws.Range("L132").Value = 14854.8
$ws.Range("M132").Value = -3459.3695
$ws.Range("N132").Value = -19914.8

$ws.Range("H136").Value = 685.0208
$ws.Range("I136").Value = 394.07693
$ws.Range("J136").Value = 1945.7778
$ws.Range("K136").Value = 1182.23079
$ws.Range("L136").Value = 5837.3334
$ws.Range("M136").Value = 1367.76921
$ws.Range("N136").Value = -10937.3334

Write-Output "edits applied"
